$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '96.325.02'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.688.16'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.39%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '236.07'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -2.72%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.89'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +4.99%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '651.27'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.71%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.424'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.60%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.05'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -1.25%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '3.685.34'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '44.32'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.204'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000291'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +12.68%  '
$ws.Range('B15').Value = 'Toncoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.70'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +2.76%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.373.82'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.50%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '96.047.24'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.690.91'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.93'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -1.63%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.77'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.72%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.75'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -10.36%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.503'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -4.63%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '518.14'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.36'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -2.61%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.0000204'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.30%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.97'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.16%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '100.81'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '13.13'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.169'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.94%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.01'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '12.05'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.96%  '
$ws.Range('B32').Value = 'Dai'
$ws.Range('C32').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('B33').Value = 'Cronos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.183'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.73%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.82'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +5.12%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '32.17'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -3.42%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '651.63'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +5.64%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.585'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.75'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.31%  '
$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.81'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +10.91%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '40.97'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -4.98%  '
$ws.Range('B43').Value = 'ImmutableX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.04'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +5.05%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.159'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +2.08%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.959'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0447'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +1.70%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.429'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +7.11%  '
$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '23.56'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.26'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -1.65%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.42'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('B51').Value = 'MantraDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.55'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +2.98%  '
